# Auto-generated Excel COM-interop edit script
# Applies targeted cell-value corrections per the commit diff
# ("reworked dur calc, lots of other changes") across the
# site_metrics, mk_duration and mk_intra_annual sheets.

$wb = $excel.ActiveWorkbook

# --- site_metrics ---
$ws = $wb.Worksheets.Item("site_metrics")
$ws.Range("O3").Value = 0.006638589657142382
$ws.Range("O4").Value = 0.01689687697811004
$ws.Range("O11").Value = 0.3414642721151773
$ws.Range("O13").Value = 0.00446236154459127
$ws.Range("AK13").Value = $true
$ws.Range("O14").Value = 0.008533913737589298
$ws.Range("AK14").Value = $true
$ws.Range("O20").Value = 0.00396137615532286
$ws.Range("O21").Value = 0.09449115251611298
$ws.Range("O34").Value = 0.01688877555810498
$ws.Range("O36").Value = 0.0009109086966998175
$ws.Range("AK36").Value = $true
$ws.Range("AK37").Value = $true
$ws.Range("O39").Value = 0.353867549199511
$ws.Range("AK40").Value = $true
$ws.Range("AK41").Value = $true
$ws.Range("O42").Value = 0.5536757223203924
$ws.Range("O43").Value = 0.003720445005879518
$ws.Range("AK43").Value = $true
$ws.Range("AK44").Value = $true
$ws.Range("O52").Value = 0.0007259547207608328
$ws.Range("AK52").Value = $true
$ws.Range("O55").Value = 0.06873677721530226
$ws.Range("O56").Value = 0.08879613705747835
$ws.Range("O58").Value = 0.7306720040936679
$ws.Range("AK60").Value = $true
$ws.Range("O62").Value = 0.0009916035592921026
$ws.Range("O67").Value = 0.01099562600072403
$ws.Range("AK69").Value = $true
$ws.Range("O72").Value = 0.007682911412505886
$ws.Range("AK73").Value = $true
$ws.Range("O76").Value = 0.0009276591751838549
$ws.Range("O79").Value = 0.01676367222418533

# --- mk_duration ---
$ws = $wb.Worksheets.Item("mk_duration")
$ws.Range("K5").Value = "no trend"
$ws.Range("L5").Value = $false
$ws.Range("M5").Value = 0.1544237061170961
$ws.Range("N5").Value = 1.424078649513432
$ws.Range("O5").Value = 0.3888888888888889
$ws.Range("P5").Value = 14
$ws.Range("Q5").Value = 83.33333333333333
$ws.Range("R5").Value = 0.9166666666666667
$ws.Range("S5").Value = -2.666666666666667
$ws.Range("M31").Value = 0.9311749148236654
$ws.Range("N31").Value = -0.08636670341750609
$ws.Range("O31").Value = -0.01231527093596059
$ws.Range("P31").Value = -5
$ws.Range("Q31").Value = 2145
$ws.Range("K61").Value = "no trend"
$ws.Range("L61").Value = $false
$ws.Range("M61").Value = 0.213590891281481
$ws.Range("N61").Value = 1.243751891458663
$ws.Range("O61").Value = 0.1897233201581028
$ws.Range("P61").Value = 48
$ws.Range("Q61").Value = 1428
$ws.Range("R61").Value = 0.08333333333333337
$ws.Range("S61").Value = 1.75
$ws.Range("M78").Value = 0.5807172923623485
$ws.Range("N78").Value = -0.5523372814706976
$ws.Range("O78").Value = -0.08
$ws.Range("P78").Value = -26
$ws.Range("Q78").Value = 2048.666666666667
$ws.Range("R78").Value = -0.02083333333333333
$ws.Range("S78").Value = 3.59375
$ws.Range("K82").Value = "no trend"
$ws.Range("L82").Value = $false
$ws.Range("M82").Value = 0.3513050661832442
$ws.Range("N82").Value = -0.9320608751309648
$ws.Range("O82").Value = -0.1978021978021978
$ws.Range("P82").Value = -18
$ws.Range("Q82").Value = 332.6666666666667
$ws.Range("R82").Value = -0.11875
$ws.Range("S82").Value = 3.549652777777778

# --- mk_intra_annual ---
$ws = $wb.Worksheets.Item("mk_intra_annual")
$ws.Range("M5").Value = 0.5084542305885602
$ws.Range("N5").Value = 0.6612465225335806
$ws.Range("O5").Value = 0.1944444444444444
$ws.Range("P5").Value = 7
$ws.Range("Q5").Value = 82.33333333333333
$ws.Range("R5").Value = 0.08333333333333333
$ws.Range("S5").Value = 0.6666666666666667
$ws.Range("M31").Value = 0.9308309497786535
$ws.Range("N31").Value = 0.08679941859050837
$ws.Range("O31").Value = 0.01231527093596059
$ws.Range("P31").Value = 5
$ws.Range("Q31").Value = 2123.666666666667
$ws.Range("M61").Value = 0.5742632659400388
$ws.Range("N61").Value = -0.5617838953175571
$ws.Range("O61").Value = -0.08695652173913043
$ws.Range("P61").Value = -22
$ws.Range("Q61").Value = 1397.333333333333
$ws.Range("S61").Value = 5
$ws.Range("M78").Value = 1
$ws.Range("N78").Value = 0
$ws.Range("O78").Value = 0.003076923076923077
$ws.Range("P78").Value = 1
$ws.Range("Q78").Value = 2015
$ws.Range("R78").Value = 0
$ws.Range("S78").Value = 4
$ws.Range("K82").Value = "no trend"
$ws.Range("L82").Value = $false
$ws.Range("M82").Value = 0.866677633709956
$ws.Range("N82").Value = 0.1678800645554932
$ws.Range("O82").Value = 0.04395604395604396
$ws.Range("P82").Value = 4
$ws.Range("Q82").Value = 319.3333333333333
$ws.Range("R82").Value = 0
$ws.Range("S82").Value = 5

